$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.580.19'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '2.027.42'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'255.72"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.68%  '
$ws.Range('D6').Value = "'0.619"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'56.83"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -8.50%  '
$ws.Range('D9').Value = "'0.383"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.66%  '
$ws.Range('D10').Value = "'0.0784"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('E11').Value = '  -2.23%  '
$ws.Range('D12').Value = "'14.52"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.80%  '
$ws.Range('D13').Value = '2.325.57'
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('D14').Value = "'0.816"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.37%  '
$ws.Range('D15').Value = "'21.09"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.28%  '
$ws.Range('D16').Value = "'5.35"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').Value = '2.028.72'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').Value = '37.519.24'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('D19').Value = "'69.58"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('D20').Value = '0.0₃0848'
$ws.Range('E20').Value = '  -2.42%  '
$ws.Range('D21').Value = "'5.20"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').Value = "'228.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('E23').Value = '  +3.54%  '
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').Value = "'2.33"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').Value = "'163.89"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'9.04"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.45%  '
$ws.Range('D28').Value = "'19.82"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('E29').Value = '  -10.46%  '
$ws.Range('D30').Value = "'1.37"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').Value = "'0.120"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('D32').Value = "'0.0666"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.25%  '
$ws.Range('D33').Value = "'4.70"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.04%  '
$ws.Range('D34').Value = "'4.57"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').Value = "'2.43"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.18%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = "'1.82"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.06%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').Value = "'3.39"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('E39').Value = '  -3.18%  '
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('D41').Value = "'0.0966"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('E42').Value = '  +1.06%  '
$ws.Range('D43').Value = "'0.0215"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.27%  '
$ws.Range('D44').Value = '1.409.13'
$ws.Range('E44').Value = '  +1.69%  '
$ws.Range('D45').Value = "'15.99"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.81%  '
$ws.Range('D46').Value = "'90.82"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').Value = "'7.32"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('D50').Value = "'2.01"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.07%  '
$ws.Range('D51').Value = '2.216.17'
$ws.Range('E51').Value = '  +1.40%  '
